$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the project name value in B3 (shared string "Bookstore Management")
$ws.Range("B3").Value = "Bookstore Management"

# Make A3 ("Project" label) bold to match header style
$ws.Range("A3").Font.Bold = $true

# Update the active selection to A8 to match final workbook state
$ws.Range("A8").Select()
